# The source data added one more weekly record at the top of the price
# history (row 9), pushing every existing record down by one row.
# Re-create that by inserting a new row at row 9 and filling it with the
# new observation; every row that used to be at N (N = 9..100) now lives
# at N+1 automatically once the row is inserted, which is exactly what the
# diff shows (each row's old values reappear one row further down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 9; rows 9-100 shift down to 10-101.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new weekly record.
$ws.Cells.Item(9, 1).Value2 = 4
$ws.Cells.Item(9, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(9, 3).Value2 = "Los Lagos"
$ws.Cells.Item(9, 4).Value2 = 44817
$ws.Cells.Item(9, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(9, 5).Value2 = 10
$ws.Cells.Item(9, 6).Value2 = 100112026
$ws.Cells.Item(9, 7).Value2 = "Haba"
$ws.Cells.Item(9, 8).Value2 = "Sin especificar"
$ws.Cells.Item(9, 9).Value2 = "Primera"
$ws.Cells.Item(9, 10).Value2 = 90
$ws.Cells.Item(9, 11).Value2 = 15000
$ws.Cells.Item(9, 12).Value2 = 15000
$ws.Cells.Item(9, 13).Value2 = 15000
$ws.Cells.Item(9, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(9, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(9, 16).Value2 = 600
$ws.Cells.Item(9, 17).Value2 = 25
$ws.Cells.Item(9, 18).Value2 = "Hortaliza"
